$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "USAPIENS" block (rows 95-154) currently labels its rows in column E
# with bare city names ("Bogotá" / "Medellín" / "Palmira"). Relabel them as
# "Puesto <Ciudad>" to match the new x-axis labels used for the chart.
for ($r = 95; $r -le 114; $r++) {
    $ws.Range("E$r").Value = "Puesto Bogotá"
}
for ($r = 115; $r -le 134; $r++) {
    $ws.Range("E$r").Value = "Puesto Medellín"
}
for ($r = 135; $r -le 154; $r++) {
    $ws.Range("E$r").Value = "Puesto Palmira"
}

# Column E now holds longer labels, so widen it to fit them.
$ws.Columns.Item(5).ColumnWidth = 31

# Match the author's final selection/view state.
$ws.Range("I149").Select()
